$d = $word.ActiveDocument

# The body is a single paragraph whose lines are separated by manual line
# breaks (w:br, which Word's text model represents as Chr(11)). Line 1
# ("1. Jane Fam tags: ...") is left completely untouched by the edit, so we
# build the Find/Replace strings starting at line 2 to avoid disturbing it.

$brk = [char]11

$oldText = (
    "2. Jonathan Dee tags: [Chemistry] [Physics] " + $brk +
    "3. Fifth " + $brk +
    "4. Mary Laking tags: [Math] " + $brk +
    "5. GetTestThree " + $brk +
    "6. GetTestTwo " + $brk +
    "7. GetTest " + $brk +
    "8. MissingGroupIDAdd " + $brk +
    "9. RemoveFailure " + $brk +
    "10. RemoveTest " + $brk +
    "11. NotFoundGroup " + $brk +
    "12. MissingGroupIDRemove " + $brk +
    "13. ExportTest "
)

$newText = (
    "2. Fifth " + $brk +
    "3. Mary Laking tags: [Math] " + $brk +
    "4. Jonathan Dee tags: [Chemistry] [Physics] " + $brk +
    "5. GetTestThree " + $brk +
    "6. GetTest " + $brk +
    "7. GetTestTwo " + $brk +
    "8. MissingGroupIDAdd " + $brk +
    "9. RemoveFailure " + $brk +
    "10. NotFoundGroup " + $brk +
    "11. RemoveTest " + $brk +
    "12. MissingGroupIDRemove " + $brk +
    "13. ExportTestTwo " + $brk +
    "14. ExportTestThree " + $brk +
    "15. ExportTest "
)

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
